$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for all existing data rows (2-349)
# from 45190 to 45192.
$ws.Range("C2:C349").Value = 45192

# Row 349 gains an explicit row height (15pt) as part of the edit.
$ws.Rows.Item(349).RowHeight = 15

# Append new row 350 with the new record.
$ws.Cells.Item(350, 1).Value = "A 44814-2023"

$ws.Range("B350").NumberFormat = "YYYY-MM-DD"
$ws.Range("B350").Value = 45190

$ws.Range("C350").NumberFormat = "YYYY-MM-DD"
$ws.Range("C350").Value = 45192

$ws.Cells.Item(350, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(350, 5).Value = "KÖPING"

$ws.Cells.Item(350, 7).Value = 24.2
$ws.Cells.Item(350, 8).Value = 0
$ws.Cells.Item(350, 9).Value = 0
$ws.Cells.Item(350, 10).Value = 0
$ws.Cells.Item(350, 11).Value = 0
$ws.Cells.Item(350, 12).Value = 0
$ws.Cells.Item(350, 13).Value = 0
$ws.Cells.Item(350, 14).Value = 0
$ws.Cells.Item(350, 15).Value = 0
$ws.Cells.Item(350, 16).Value = 0
$ws.Cells.Item(350, 17).Value = 0

$ws.Range("R350").WrapText = $true
$ws.Range("R350").Value = ""
